$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the body/ellipsis text for the affected D-column cells
# ("…" -> "...") as part of the editorial text cleanup.
$ws.Cells.Item(3, 4).Value2 = 'T148_49 - Po przelicz. KP, kwota ponownie oblicz. i zwal. KP - na Pana koncie wynosi ...'
$ws.Cells.Item(4, 4).Value2 = 'T16_66 - Art. 55a - ... prawo do wcześn. em. zostanie przywrócone w przyp. wycofania wniosku ...'
$ws.Cells.Item(5, 4).Value2 = 'T2_151 - Art. 26a - W związku z przyznaniem emerytury w powsz. wieku em. ... prawo do em. wcześn. ustaje z dniem ...'
$ws.Cells.Item(6, 4).Value2 = 'T2_152 - Art. 26a - W związku z obliczeniem emerytury w powsz. wieku em. ... prawo do em. wcześn. ustaje z dniem ...'
$ws.Cells.Item(8, 4).Value2 = 'T8_159 - Art 174 ust. 2a i 185a - Po ponownym ustaleniu wysokość ... wynosi: SKL + KAP / ŚDTŻ'
$ws.Cells.Item(9, 4).Value2 = 'T8_160 - Art 174 ust. 2a i 185a - Po ponownym ustaleniu wysokość ... wynosi: KAP / ŚDTZ'
$ws.Cells.Item(18, 4).Value2 = 'T8_169 - Art. 26 ust. 5 i 6 - ŚDTŻ ... w dniu spełnienia warunków'
$ws.Cells.Item(21, 4).Value2 = 'T88_147 - Art. 55a - Podstawę obliczenia emerytury stanowi ... kwota składek + KP'
$ws.Cells.Item(22, 4).Value2 = 'T88_148 - Art. 55a - Podstawę obliczenia emerytury stanowi ... kwota składek'
$ws.Cells.Item(23, 4).Value2 = 'T88_149 - Art. 55a - Podstawę obliczenia emerytury stanowi ... KP'
$ws.Cells.Item(25, 4).Value2 = 'T148_49 - Po przelicz. KP, kwota ponownie oblicz. i zwal. KP - na Pana koncie wynosi ...'
$ws.Cells.Item(27, 4).Value2 = 'T8_159 - Art 174 ust. 2a i 185a - Po ponownym ustaleniu wysokość ... wynosi: SKL + KAP / ŚDTŻ'
$ws.Cells.Item(28, 4).Value2 = 'T8_160 - Art 174 ust. 2a i 185a - Po ponownym ustaleniu wysokość ... wynosi: KAP / ŚDTZ'
$ws.Cells.Item(30, 4).Value2 = 'T148_49 - Po przelicz. KP, kwota ponownie oblicz. i zwal. KP - na Pana koncie wynosi ...'
$ws.Cells.Item(32, 4).Value2 = 'T8_159 - Art 174 ust. 2a i 185a - Po ponownym ustaleniu wysokość ... wynosi: SKL + KAP / ŚDTŻ'
$ws.Cells.Item(33, 4).Value2 = 'T8_160 - Art 174 ust. 2a i 185a - Po ponownym ustaleniu wysokość ... wynosi: KAP / ŚDTZ'
$ws.Cells.Item(35, 4).Value2 = 'T148_49 - Po przelicz. KP, kwota ponownie oblicz. i zwal. KP - na Pana koncie wynosi ...'
$ws.Cells.Item(37, 4).Value2 = 'T8_159 - Art 174 ust. 2a i 185a - Po ponownym ustaleniu wysokość ... wynosi: SKL + KAP / ŚDTŻ'
$ws.Cells.Item(38, 4).Value2 = 'T8_160 - Art 174 ust. 2a i 185a - Po ponownym ustaleniu wysokość ... wynosi: KAP / ŚDTZ'
$ws.Cells.Item(47, 4).Value2 = 'T8_169 - Art. 26 ust. 5 i 6 - ŚDTŻ ... w dniu spełnienia warunków'
$ws.Cells.Item(51, 4).Value2 = 'T148_49 - Po przelicz. KP, kwota ponownie oblicz. i zwal. KP - na Pana koncie wynosi ...'
$ws.Cells.Item(53, 4).Value2 = 'T8_159 - Art 174 ust. 2a i 185a - Po ponownym ustaleniu wysokość ... wynosi: SKL + KAP / ŚDTŻ'
$ws.Cells.Item(54, 4).Value2 = 'T8_160 - Art 174 ust. 2a i 185a - Po ponownym ustaleniu wysokość ... wynosi: KAP / ŚDTZ'
$ws.Cells.Item(58, 4).Value2 = 'T10b_66 - Art. 26 ust. 5 i 6 - Do ustalenia wys. okr. em. przyjęto ŚDTŻ ... w dniu zgłoszenia wniosku o em. ... korzystniejsze od ... w dniu osiągn. powsz. w. em. '
$ws.Cells.Item(59, 4).Value2 = 'T10b_67 - Art. 26 ust. 5 i 6 - Do ustalenia wys. okr. em. przyjęto ŚDTŻ ... w dniu osiągn. powsz. w. em. ... korzystniejsze od ... w dniu zgłoszenia wniosku o em.'
$ws.Cells.Item(60, 4).Value2 = 'T10b_68 - Art. 26 ust. 5 i 6 - Do ustalenia wys. okr. em. przyjęto ŚDTŻ ... w dniu osiągn. powsz. w. em. ... korzystniejsze od ... wypłata zawiesz. renty z tyt. niezd. do pracy'
$ws.Cells.Item(61, 4).Value2 = 'T10b_69 - Art. 26 ust. 5 i 6 - Do ustalenia wys. okr. em. przyjęto ŚDTŻ ... wypłata zawiesz. renty z tyt. niezd. do pracy ... korzystniejsze od ... osiągn. powsz. w. em.'
$ws.Cells.Item(67, 4).Value2 = 'T148_50 - Po przelicz. KP, kwota ponownie oblicz. i zwal. KP - na koncie osoby zmarłej wynosi ...'
$ws.Cells.Item(70, 4).Value2 = 'T8_172 - Art 174 ust. 2a i 185a - Po ponownym ustaleniu ... wysokość świadcz. osoby zmarłej ... wynosi: SKL + KAP / ŚDTŻ'
$ws.Cells.Item(71, 4).Value2 = 'T8_173 - Art 174 ust. 2a i 185a - Po ponownym ustaleniu ... wysokość świadcz. osoby zmarłej ... wynosi: KAP / ŚDTŻ'
$ws.Cells.Item(74, 4).Value2 = 'T8_176 - Art. 55 i 55a - Po ponownym ustaleniu ... wysokość świadcz. osoby zmarłej ... wynosi: SKL + KAP - POBR.EM. / ŚDTŻ'
$ws.Cells.Item(75, 4).Value2 = 'T8_177 - Art. 55 i 55a - Po ponownym ustaleniu ... wysokość świadcz. osoby zmarłej ... wynosi: SKL - POBR.EM. / ŚDTŻ'
$ws.Cells.Item(76, 4).Value2 = 'T8_178 - Art. 55 i 55a - Po ponownym ustaleniu ... wysokość świadcz. osoby zmarłej ... wynosi: KAP - POBR.EM. / ŚDTŻ'
$ws.Cells.Item(81, 4).Value2 = 'T8_183 - Art. 26 ust. 5 i 6 - ŚDTŻ ... w dniu spełnienia warunków - dla osoby zmarłej'
$ws.Cells.Item(86, 4).Value2 = 'T88_150 - Art. 55a - Podstawę obliczenia emerytury osoby zmarłej stanowi ... kwota składek + KP'
$ws.Cells.Item(87, 4).Value2 = 'T88_151 - Art. 55a - Podstawę obliczenia emerytury osoby zmarłej stanowi ... kwota składek'
$ws.Cells.Item(88, 4).Value2 = 'T88_152 - Art. 55a - Podstawę obliczenia emerytury osoby zmarłej stanowi ... KP'
$ws.Cells.Item(96, 4).Value2 = 'T8_183 - Art. 26 ust. 5 i 6 - ŚDTŻ ... w dniu spełnienia warunków - dla osoby zmarłej'
$ws.Cells.Item(101, 4).Value2 = 'T16_69 - Art. 110a - Ponownego ustalenia em. osoby zmarłej dokonano ... może nastąpić wyłącznie jeden raz.'
$ws.Cells.Item(103, 4).Value2 = 'T148_49 - Po przelicz. KP, kwota ponownie oblicz. i zwal. KP - na Pana koncie wynosi ...'
$ws.Cells.Item(105, 4).Value2 = 'T8_159 - Art 174 ust. 2a i 185a - Po ponownym ustaleniu wysokość ... wynosi: SKL + KAP / ŚDTŻ'
$ws.Cells.Item(106, 4).Value2 = 'T8_160 - Art 174 ust. 2a i 185a - Po ponownym ustaleniu wysokość ... wynosi: KAP / ŚDTZ'
$ws.Cells.Item(109, 4).Value2 = 'T8_163 - Art. 26 ust. 5 i 6 - z urzędu ŚDTŻ ... w dniu osiągn. powsz. w. em. ... korzystniejsze od ... wypłata zawiesz. renty z tyt. niezd. do pracy'
$ws.Cells.Item(115, 4).Value2 = 'T8_169 - Art. 26 ust. 5 i 6 - ŚDTŻ ... w dniu spełnienia warunków'
$ws.Cells.Item(119, 4).Value2 = 'T148_49 - Po przelicz. KP, kwota ponownie oblicz. i zwal. KP - na Pana koncie wynosi ...'
$ws.Cells.Item(121, 4).Value2 = 'T8_159 - Art 174 ust. 2a i 185a - Po ponownym ustaleniu wysokość ... wynosi: SKL + KAP / ŚDTŻ'
$ws.Cells.Item(122, 4).Value2 = 'T8_160 - Art 174 ust. 2a i 185a - Po ponownym ustaleniu wysokość ... wynosi: KAP / ŚDTZ'
$ws.Cells.Item(134, 4).Value2 = 'T8_183 - Art. 26 ust. 5 i 6 - ŚDTŻ ... w dniu spełnienia warunków - dla osoby zmarłej'
$ws.Cells.Item(139, 4).Value2 = 'T16_69 - Art. 110a - Ponownego ustalenia em. osoby zmarłej dokonano ... może nastąpić wyłącznie jeden raz.'
$ws.Cells.Item(141, 4).Value2 = 'T50_449 - Zgodnie z art. 174 ust. 2a ... przy ustalaniu KP okresy urlopów lub niewyk. pracy oblicza się przyj. 1,3% PW ...'
$ws.Cells.Item(142, 4).Value2 = 'T50_450 - Zgodnie z art. 185a ... poprzez dodanie do okresów nieskładk. okresów studiów wyższ. w wymiarze 1/3 okr. składk.'
$ws.Cells.Item(143, 4).Value2 = 'T50_451 - Art. 55a - Wprowadzony ustawą ... umożliwia ponowne obliczenie, z zast. art. 55 ... osobie, która miała ustalone prawo do em. wcześn. i kontynuowała ubezp.'
$ws.Cells.Item(144, 4).Value2 = 'T50_452 - Art. 26 ust. 6 - Wysokość em. podlega ponownemu ustaleniu ... jeżeli obowiązywała inna tablica ŚDTŻ ...'
$ws.Cells.Item(147, 4).Value2 = 'T54_698 - Art. 55 i 55a - nie kontynuował/a Pan/i ubezp. ... po osiągn. powsz. w. em. ... brak podstaw do oblicz.'
$ws.Cells.Item(148, 4).Value2 = 'T54_699 - Art. 55 i 55a - emerytura została przyznana przed dniem 1 stycznia 2009 r., a zatem brak jest podstaw do obliczenia ...'
$ws.Cells.Item(150, 4).Value2 = 'T54_701 - wniosek o emeryturę zgłosił Pan(i) w miesiącu, w którym obowiązywała ta sama tablica ŚDTŻ ...'
$ws.Cells.Item(151, 4).Value2 = 'T50_453 - Zgodnie z art. 110a ... dotyczy wyłącznie em. i może nastąpić wyłącznie jeden raz'
$ws.Cells.Item(155, 4).Value2 = 'T54_702 - Art. 110a - do ponownego ustalenia ... wskazano wynagrodzenia przyjęte do ostatnio obliczonej podstawy wymiaru'
$ws.Cells.Item(156, 4).Value2 = 'T54_703 - Art. 110a - do ponownego ustalenia podstawy wymiaru emerytury nie wskazano wynagrodzeń ...'
$ws.Cells.Item(157, 4).Value2 = 'T54_704 - Art. 110a - nie podlegał/a Pan/i ubezp. społecz. ... po przyznaniu em. ... brak podstaw do ponownego przel.'
$ws.Cells.Item(158, 4).Value2 = 'T50_454 - Art. 110a - Wysokość emerytury podlega ponownemu ustaleniu ... gdy zostały spełnione warunki ...'
$ws.Cells.Item(160, 4).Value2 = 'T50_449 - Zgodnie z art. 174 ust. 2a ... przy ustalaniu KP okresy urlopów lub niewyk. pracy oblicza się przyj. 1,3% PW ...'
$ws.Cells.Item(161, 4).Value2 = 'T50_450 - Zgodnie z art. 185a ... poprzez dodanie do okresów nieskładk. okresów studiów wyższ. w wymiarze 1/3 okr. składk.'
$ws.Cells.Item(165, 4).Value2 = 'T50_449 - Zgodnie z art. 174 ust. 2a ... przy ustalaniu KP okresy urlopów lub niewyk. pracy oblicza się przyj. 1,3% PW ...'
$ws.Cells.Item(166, 4).Value2 = 'T50_450 - Zgodnie z art. 185a ... poprzez dodanie do okresów nieskładk. okresów studiów wyższ. w wymiarze 1/3 okr. składk.'
$ws.Cells.Item(169, 4).Value2 = 'T148_49 - Po przelicz. KP, kwota ponownie oblicz. i zwal. KP - na Pana koncie wynosi ...'
$ws.Cells.Item(170, 4).Value2 = 'T8_159 - Art 174 ust. 2a i 185a - Po ponownym ustaleniu wysokość ... wynosi: SKL + KAP / ŚDTŻ'
$ws.Cells.Item(171, 4).Value2 = 'T8_160 - Art 174 ust. 2a i 185a - Po ponownym ustaleniu wysokość ... wynosi: KAP / ŚDTZ'

# Window + selection cosmetic changes
$excel.ActiveWindow.Left = 3750
$ws.Range("D3").Select()
